# Auto-generated Excel COM-interop script to apply numeric value updates
# as described in the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 19843
$ws.Range("J7").Value = 24052.5
$ws.Range("L7").Value = 24052.5
$ws.Range("N7").Value = -24276.5

$ws.Range("H14").Value = 19843
$ws.Range("J14").Value = 24052.5
$ws.Range("L14").Value = 24052.5
$ws.Range("N14").Value = -24434.5

$ws.Range("H40").Value = 1741.4078
$ws.Range("I40").Value = 1709.0161
$ws.Range("J40").Value = 1884.8572
$ws.Range("K40").Value = 1709.0161
$ws.Range("L40").Value = 1884.8572
$ws.Range("M40").Value = -1534.0161
$ws.Range("N40").Value = -2234.8572

$ws.Range("H43").Value = 1185.3636
$ws.Range("I43").Value = 866.6667
$ws.Range("J43").Value = 1304.875
$ws.Range("K43").Value = 866.6667
$ws.Range("L43").Value = 1304.875
$ws.Range("M43").Value = -797.6667
$ws.Range("N43").Value = -1442.875

$ws.Range("H80").Value = 608.2059
$ws.Range("I80").Value = 569.6
$ws.Range("J80").Value = 663.3570999999999
$ws.Range("K80").Value = 1708.8
$ws.Range("L80").Value = 1990.0713
$ws.Range("M80").Value = -710.8000000000002
$ws.Range("N80").Value = -3986.0713

$ws.Range("H83").Value = 608.2059
$ws.Range("I83").Value = 569.6
$ws.Range("J83").Value = 663.3570999999999
$ws.Range("K83").Value = 5126.400000000001
$ws.Range("L83").Value = 5970.2139
$ws.Range("M83").Value = -134.4000000000005
$ws.Range("N83").Value = -15954.2139

$ws.Range("H107").Value = 17857392
$ws.Range("I107").Value = 17857392
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 17857392
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -17855472
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5088.078
$ws.Range("I32").Value = 4513.25
$ws.Range("K32").Value = 4513.25
$ws.Range("M32").Value = -4226.25

$ws.Range("H102").Value = 3706140.5
$ws.Range("I102").Value = 4117600.5
$ws.Range("K102").Value = 4117600.5
$ws.Range("M102").Value = -4115978.5

$ws.Range("H132").Value = 2331.9788
$ws.Range("I132").Value = 1016.0714
$ws.Range("J132").Value = 4271.2104
$ws.Range("K132").Value = 3048.2142
$ws.Range("L132").Value = 12813.6312
$ws.Range("M132").Value = -518.2142000000003
$ws.Range("N132").Value = -17873.6312

$ws.Range("H135").Value = 57841.375
$ws.Range("J135").Value = 57841.375
$ws.Range("L135").Value = 57841.375
$ws.Range("N135").Value = -67981.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9599.423000000001
$ws.Range("I105").Value = 12553.389
$ws.Range("K105").Value = 12553.389
$ws.Range("M105").Value = -10806.389

$ws.Range("H107").Value = 1086.9231
$ws.Range("I107").Value = 1093.0834
$ws.Range("K107").Value = 1093.0834
$ws.Range("M107").Value = 826.9166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6172.206
$ws.Range("I31").Value = 1575.3158
$ws.Range("J31").Value = 11994.934
$ws.Range("K31").Value = 1575.3158
$ws.Range("L31").Value = 11994.934
$ws.Range("M31").Value = -1280.3158
$ws.Range("N31").Value = -12584.934

$ws.Range("H34").Value = 6172.206
$ws.Range("I34").Value = 1575.3158
$ws.Range("J34").Value = 11994.934
$ws.Range("K34").Value = 1575.3158
$ws.Range("L34").Value = 11994.934
$ws.Range("M34").Value = -1373.3158
$ws.Range("N34").Value = -12398.934

$ws.Range("H58").Value = 1650.4828
$ws.Range("I58").Value = 1203.125
$ws.Range("J58").Value = 2201.077
$ws.Range("K58").Value = 1203.125
$ws.Range("L58").Value = 2201.077
$ws.Range("M58").Value = -1000.125
$ws.Range("N58").Value = -2607.077

$ws.Range("H136").Value = 1650.4828
$ws.Range("I136").Value = 1203.125
$ws.Range("J136").Value = 2201.077
$ws.Range("K136").Value = 3609.375
$ws.Range("L136").Value = 6603.231000000001
$ws.Range("M136").Value = -1059.375
$ws.Range("N136").Value = -11703.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 300307.75
$ws.Range("I5").Value = 354.58334
$ws.Range("K5").Value = 1063.75002
$ws.Range("M5").Value = -951.7500199999999

$ws.Range("H68").Value = 2019
$ws.Range("I68").Value = 422.8
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1268.4
$ws.Range("L68").Value = 30000
$ws.Range("M68").Value = -457.4000000000001
$ws.Range("N68").Value = -31622

$ws.Range("H71").Value = 2019
$ws.Range("I71").Value = 422.8
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 3805.2
$ws.Range("L71").Value = 90000
$ws.Range("M71").Value = 250.7999999999997
$ws.Range("N71").Value = -98112

$ws.Range("H97").Value = 10000686
$ws.Range("I97").Value = 25000376
$ws.Range("J97").Value = 894
$ws.Range("K97").Value = 75001128
$ws.Range("L97").Value = 2682
$ws.Range("M97").Value = -75000632
$ws.Range("N97").Value = -3674

$ws.Range("H129").Value = 1823.3334
$ws.Range("I129").Value = 801.6667
$ws.Range("J129").Value = 3866.6667
$ws.Range("K129").Value = 2405.0001
$ws.Range("L129").Value = 11600.0001
$ws.Range("M129").Value = 2594.9999
$ws.Range("N129").Value = -21600.0001

$ws.Range("H131").Value = 1352305.2
$ws.Range("I131").Value = 6667334.5
$ws.Range("J131").Value = 1026.5593
$ws.Range("K131").Value = 20002003.5
$ws.Range("L131").Value = 3079.6779
$ws.Range("M131").Value = -19996963.5
$ws.Range("N131").Value = -13159.6779

$ws.Range("H132").Value = 7939999.5
$ws.Range("I132").Value = 970
$ws.Range("J132").Value = 11115612
$ws.Range("K132").Value = 8730
$ws.Range("L132").Value = 100040508
$ws.Range("M132").Value = -6200
$ws.Range("N132").Value = -100045568

$ws.Range("H135").Value = 300307.75
$ws.Range("I135").Value = 354.58334
$ws.Range("K135").Value = 3191.25006
$ws.Range("M135").Value = -656.2500600000003

$ws.Range("H137").Value = 12838551
$ws.Range("I137").Value = 13028.777
$ws.Range("J137").Value = 19628534
$ws.Range("K137").Value = 39086.331
$ws.Range("L137").Value = 58885602
$ws.Range("M137").Value = -33986.331
$ws.Range("N137").Value = -58895802

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10753517
$ws.Range("I46").Value = 14493384
$ws.Range("J46").Value = 1399.75
$ws.Range("K46").Value = 14493384
$ws.Range("L46").Value = 1399.75
$ws.Range("M46").Value = -14493196
$ws.Range("N46").Value = -1775.75

$ws.Range("H55").Value = 93750280
$ws.Range("I55").Value = 100000300
$ws.Range("J55").Value = 83333576
$ws.Range("K55").Value = 100000300
$ws.Range("L55").Value = 83333576
$ws.Range("M55").Value = -100000127
$ws.Range("N55").Value = -83333922

$ws.Range("H68").Value = 35716228
$ws.Range("I68").Value = 1498.25
$ws.Range("J68").Value = 83335864
$ws.Range("K68").Value = 1498.25
$ws.Range("L68").Value = 83335864
$ws.Range("M68").Value = -749.25
$ws.Range("N68").Value = -83337362

$ws.Range("H71").Value = 35716228
$ws.Range("I71").Value = 1498.25
$ws.Range("J71").Value = 83335864
$ws.Range("K71").Value = 7491.25
$ws.Range("L71").Value = 416679320
$ws.Range("M71").Value = -3747.25
$ws.Range("N71").Value = -416686808

$ws.Range("H122").Value = 8931197
$ws.Range("I122").Value = 71428570
$ws.Range("J122").Value = 3000.7144
$ws.Range("K122").Value = 214285710
$ws.Range("L122").Value = 9002.143199999999
$ws.Range("M122").Value = -214283260
$ws.Range("N122").Value = -13902.1432

$ws.Range("H132").Value = 14328902
$ws.Range("I132").Value = 18339842
$ws.Range("J132").Value = 4114.7144
$ws.Range("K132").Value = 55019526
$ws.Range("L132").Value = 12344.1432
$ws.Range("M132").Value = -55016996
$ws.Range("N132").Value = -17404.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 1007500
$ws.Range("J11").Value = 1007500
$ws.Range("L11").Value = 1007500
$ws.Range("N11").Value = -1007784

$ws.Range("H20").Value = 500
$ws.Range("I20").Value = 500
$ws.Range("K20").Value = 500
$ws.Range("M20").Value = -260

$ws.Range("H122").Value = 2878.818
$ws.Range("I122").Value = 2349.1333
$ws.Range("K122").Value = 7047.3999
$ws.Range("M122").Value = -4597.3999

$ws.Range("H132").Value = 2113.0833
$ws.Range("I132").Value = 1600.5
$ws.Range("J132").Value = 2625.6667
$ws.Range("K132").Value = 4801.5
$ws.Range("L132").Value = 7877.000100000001
$ws.Range("M132").Value = -2271.5
$ws.Range("N132").Value = -12937.0001
